# Update the cryptos price/volume table with the latest scraped values.
# Numeric-looking "Price" (column D) values are prefixed with a leading
# apostrophe so Excel stores them as text (matching the original inline
# strings) instead of converting them to floating point numbers, which
# would introduce rounding artifacts (e.g. 135.57 -> 135.56999999999999).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.014.05'
$ws.Range("D3").Value = '2.406.89'
$ws.Range("E3").Value = '  -0.31%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''554.37'
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("D6").Value = '''135.57'
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -0.96%  '
$ws.Range("E9").Value = '  -0.27%  '
$ws.Range("D10").Value = '''5.63'
$ws.Range("E10").Value = '  -1.01%  '
$ws.Range("E11").Value = '  -0.37%  '
$ws.Range("E12").Value = '  -1.30%  '
$ws.Range("D13").Value = '''24.61'
$ws.Range("E13").Value = '  -0.19%  '
$ws.Range("D14").Value = '2.837.54'
$ws.Range("E14").Value = '  -0.23%  '
$ws.Range("D15").Value = '59.888.34'
$ws.Range("E16").Value = '  +0.53%  '
$ws.Range("D17").Value = '2.405.75'
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").Value = '''4.49'
$ws.Range("E19").Value = '  +3.23%  '
$ws.Range("D20").Value = '''326.75'
$ws.Range("E20").Value = '  -0.77%  '
$ws.Range("D21").Value = '''6.77'
$ws.Range("E21").Value = '  +1.25%  '
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '''64.62'
$ws.Range("E23").Value = '  -1.55%  '
$ws.Range("E24").Value = '  +4.54%  '
$ws.Range("D25").Value = '''8.59'
$ws.Range("E25").Value = '  +0.44%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  +4.05%  '
$ws.Range("E28").Value = '  +1.39%  '
$ws.Range("E29").Value = '  -0.96%  '
$ws.Range("D30").Value = '''169.56'
$ws.Range("E30").Value = '  -0.40%  '
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("E32").Value = '  +8.42%  '
$ws.Range("E33").Value = '  -1.97%  '
$ws.Range("D34").Value = '''18.39'
$ws.Range("E34").Value = '  -1.15%  '
$ws.Range("E35").Value = '  +0.09%  '
$ws.Range("E36").Value = '  +3.28%  '
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("E38").Value = '  +0.45%  '
$ws.Range("D39").Value = '''322.78'
$ws.Range("E39").Value = '  +2.89%  '
$ws.Range("E40").Value = '  -0.33%  '
$ws.Range("D41").Value = '''146.94'
$ws.Range("E41").Value = '  +6.40%  '
$ws.Range("E42").Value = '  -1.97%  '
$ws.Range("D43").Value = '''0.0965'
$ws.Range("E43").Value = '  -0.43%  '
$ws.Range("D44").Value = '''19.82'
$ws.Range("E44").Value = '  +2.38%  '
$ws.Range("E45").Value = '  -0.50%  '
$ws.Range("E46").Value = '  -0.43%  '
$ws.Range("E47").Value = '  -1.19%  '
$ws.Range("D48").Value = '''11.06'
$ws.Range("E48").Value = '  +0.13%  '
$ws.Range("E49").Value = '  -0.77%  '
$ws.Range("E50").Value = '  -0.66%  '
$ws.Range("E51").Value = '  -1.01%  '
